# "finish body graphic, living spells working"
# Update the modifier picks on Sheet2 (the active sheet) and move the
# selection to reflect where the user left off (E7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Talisman"
$ws.Range("C3").Value = "Living"
$ws.Range("D3").Value = "Water"
$ws.Range("E3").Value = "Roll Count"
$ws.Range("F3").Value = "None"

$ws.Range("E7").Select()
